$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H, shifting old H (Recoup30k) and I (MeanAge)
# to I and J respectively.
$ws.Columns("H").Insert()

# Header for the new column
$ws.Range("H1").Value = "Recoup20k"

# New "Recoup20k" values for rows 2-16
$ws.Range("H2").Value = 46.1
$ws.Range("H3").Value = 29.8
$ws.Range("H4").Value = 51
$ws.Range("H5").Value = 71.1
$ws.Range("H6").Value = 66.1
$ws.Range("H7").Value = 7.17
$ws.Range("H8").Value = 13.6
$ws.Range("H9").Value = 12.5
$ws.Range("H10").Value = 9.68
$ws.Range("H11").Value = 23.6
$ws.Range("H12").Value = 53.2
$ws.Range("H13").Value = 36.2
$ws.Range("H14").Value = 61.2
$ws.Range("H15").Value = 85.7
$ws.Range("H16").Value = 79.5

# Match the author's final cursor position as seen in the saved file
$ws.Range("H17").Select()
